$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# Row 6 - shifted-in data (previously row 7's child) with refreshed path/time values
$ws.Range("B6").Value = "'6"
$ws.Range("C6").Value = "Ema " + $nbsp
$ws.Range("D6").Value = "Ardell " + $nbsp
$ws.Range("E6").Value = "-9.79,-8.09"
$ws.Range("F6").Value = "Carley(grandmother): 0533587167"
$ws.Range("H6").Value = "'37.0"

# Row 7
$ws.Range("B7").Value = "'0"
$ws.Range("C7").Value = "Trudie " + $nbsp
$ws.Range("D7").Value = "Fleta " + $nbsp
$ws.Range("E7").Value = "-4.18,-8.88"
$ws.Range("F7").Value = "Anneliese(father): 0548973345"
$ws.Range("G7").Value = "7:07:00"
$ws.Range("H7").Value = "'30.0"

# Row 8
$ws.Range("B8").Value = "'3"
$ws.Range("C8").Value = "Alexia " + $nbsp
$ws.Range("D8").Value = "Ramonita " + $nbsp
$ws.Range("E8").Value = "-1.65,-8.14"
$ws.Range("F8").Value = "Han(father): 0567537032"
$ws.Range("G8").Value = "7:11:00"
$ws.Range("H8").Value = "'26.0"

# Row 9
$ws.Range("B9").Value = "'1"
$ws.Range("C9").Value = "Corene " + $nbsp
$ws.Range("D9").Value = "Myra " + $nbsp
$ws.Range("E9").Value = "4.52,-9.26"
$ws.Range("F9").Value = "Georgie(mother): 0544823581"
$ws.Range("G9").Value = "7:20:00"
$ws.Range("H9").Value = "'17.0"

# Row 10
$ws.Range("B10").Value = "'8"
$ws.Range("C10").Value = "Marni " + $nbsp
$ws.Range("D10").Value = "Shanika " + $nbsp
$ws.Range("E10").Value = "5.4,-6.02"
$ws.Range("F10").Value = "Lady(mother): 0560804012"
$ws.Range("G10").Value = "7:24:00"
$ws.Range("H10").Value = "'13.0"

# Row 11
$ws.Range("B11").Value = "'2"
$ws.Range("C11").Value = "Elwanda " + $nbsp
$ws.Range("D11").Value = "Cassy " + $nbsp
$ws.Range("E11").Value = "-1.98,-2.1"
$ws.Range("F11").Value = "Tamisha(mother): 0550693864"
$ws.Range("G11").Value = "7:34:00"
$ws.Range("H11").Value = "'3.0"

# Row 12 (school row) - only G12 (pickup time) changes
$ws.Range("G12").Value = "7:37:00"

# Row 14 (time) - only B14 changes
$ws.Range("B14").Value = "'37.0"
